$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# --- Column B (Paterno) for rows 2-5, written first so new shared strings
#     are appended in this exact order: LOPEZ, DE LOS SANTOS, GIL, BERINSTAIN
$ws.Range("B2").Value = "LOPEZ"
$ws.Range("B3").Value = "DE LOS SANTOS"
$ws.Range("B4").Value = "GIL"
$ws.Range("B5").Value = "BERINSTAIN"

# --- Column C (Materno) for rows 2-5, written next so new shared strings
#     are appended in this exact order: ZAMUDIO, GONZALEZ, BANDALA, SAN JUAN
$ws.Range("C2").Value = "ZAMUDIO"
$ws.Range("C3").Value = "GONZALEZ"
$ws.Range("C4").Value = "BANDALA"
$ws.Range("C5").Value = "SAN JUAN"

# --- Column D (Nombres) for rows 2-5, written last so new shared strings
#     are appended in this exact order: EZRA, MARIA FERNANDA, AELEN, LUIS FERNANDO
$ws.Range("D2").Value = "EZRA"
$ws.Range("D3").Value = "MARIA FERNANDA"
$ws.Range("D4").Value = "AELEN"
$ws.Range("D5").Value = "LUIS FERNANDO"

# --- Column A (NC / student id)
$ws.Range("A2").Value = 18330051920159
$ws.Range("A3").Value = 18330051920017
$ws.Range("A4").Value = 18330051920368
$ws.Range("A5").Value = 18330051920088

# --- Column E (Nombre_Largo / subject)
$ws.Range("E2").Value = "TEMAS DE FÍSICA"
$ws.Range("E3").Value = "TEMAS DE FÍSICA"
$ws.Range("E4").Value = "TEMAS DE FÍSICA"
$ws.Range("E5").Value = "TEMAS DE ADMINISTRACIÓN"

# --- Column F (Grupo)
$ws.Range("F2").Value = "6AEM"
$ws.Range("F3").Value = "6AEV"
$ws.Range("F4").Value = "6APM"
$ws.Range("F5").Value = "6ARHV"

# --- Column G (Reprobadas)
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 1
